{"js": "// Add a new custom paragraph style \"CompactList\" (\"Compact List\"),\n// cloned from the existing \"Compact\" style: based on Body Text, quick\n// style, with 1.8pt (36 twips) spacing before/after.\n\n// Create the style (Word appends new custom styles to the style sheet).\ncontext.document.addStyle(\"Compact List\", Word.StyleType.paragraph);\nawait context.sync();\n\n// Re-fetch the freshly created style by name so property writes stick,\n// then give it the same formatting as the existing \"Compact\" style.\nconst newStyle = context.document.getStyles().getByName(\"Compact List\");\nnewStyle.baseStyle = \"BodyText\";\nnewStyle.quickStyle = true;\nnewStyle.paragraphFormat.spaceBefore = 1.8;  // 36 twips\nnewStyle.paragraphFormat.spaceAfter = 1.8;   // 36 twips\nawait context.sync();\n", "ps1": "# Add a new custom paragraph style \"CompactList\" (\"Compact List\"),\n# cloned from the existing \"Compact\" style: based on Body Text, quick\n# style, with 1.8pt (36 twips) spacing before/after.\n\n$d = $word.ActiveDocument\n\n# wdStyleTypeParagraph = 1\n$s = $d.Styles.Add(\"Compact List\", 1)\n$s.BaseStyle = $d.Styles(\"BodyText\")\n$s.QuickStyle = $true\n$s.ParagraphFormat.SpaceBefore = 1.8\n$s.ParagraphFormat.SpaceAfter = 1.8\n"}
